$d = $word.ActiveDocument

# Locate the run of text that needs to be split so the body-part name
# "veines" is wrapped in <bp>...</bp> markup tags (rendered in the
# document's blue Courier-New "tag" style, matching the other <bp> tags
# already present elsewhere in this document).
$target = $d.Content
$found = $target.Find.Execute(
    " ont les veines hemorroidales eminentes",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the target sentence to edit."
}

$startPos = $target.Start

$piece1 = " ont les "
$piece2 = "<bp>"
$piece3 = "veines"
$piece4 = "</bp>"
$piece5 = " hemorroidales eminentes"

# Replace the whole matched span with the plain concatenation first; this
# keeps a single contiguous range we can then re-slice into the five runs
# below (each slice becomes its own run once distinct formatting is
# applied to it).
$target.Text = $piece1 + $piece2 + $piece3 + $piece4 + $piece5

$pos = $startPos
$r1 = $d.Range($pos, $pos + $piece1.Length); $pos += $piece1.Length
$r2 = $d.Range($pos, $pos + $piece2.Length); $pos += $piece2.Length
$r3 = $d.Range($pos, $pos + $piece3.Length); $pos += $piece3.Length
$r4 = $d.Range($pos, $pos + $piece4.Length); $pos += $piece4.Length
$r5 = $d.Range($pos, $pos + $piece5.Length); $pos += $piece5.Length

# Plain-text runs keep the surrounding black colour.
$r1.Font.Color = 0
$r3.Font.Color = 0
$r5.Font.Color = 0

# The literal "<bp>" / "</bp>" markup runs use the small blue Courier New
# "tag" styling used throughout this document.
foreach ($tagRun in @($r2, $r4)) {
    $tagRun.Font.Name = "Courier New"
    $tagRun.Font.Color = 16711680
    $tagRun.Font.Size = 9
}
